$d = $word.ActiveDocument

# --- 1. First paragraph: pad with two trailing spaces, then append a
#        red-colored "(This is a change - Version for main branch)" note,
#        built from three separate runs (matching the authored edit). ---
$para1 = $d.Paragraphs(1).Range
$found = $para1.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

if ($found) {
    $endOfPara1 = $d.Paragraphs(1).Range
    $insertPoint = $d.Range($endOfPara1.End - 1, $endOfPara1.End - 1)
    $insertPoint.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
    $insertPoint.Font.Color = 255

    $insertPoint2 = $d.Range($d.Paragraphs(1).Range.End - 1, $d.Paragraphs(1).Range.End - 1)
    $insertPoint2.InsertAfter("rsion for main branch")
    $insertPoint2.Font.Color = 255

    $insertPoint3 = $d.Range($d.Paragraphs(1).Range.End - 1, $d.Paragraphs(1).Range.End - 1)
    $insertPoint3.InsertAfter(")")
    $insertPoint3.Font.Color = 255
}

# --- 2. Remove the trailing "ank God almighty, we are free at last."
#        paragraph that follows "Shall be lifted-nevermore!" ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
if ($lastPara.Range.Text -match "God almighty") {
    $lastPara.Range.Delete()
}
